# Remove the clip1_awake / clip2_awake / clip1_asleep / clip2_asleep /
# "poor data quality" columns (O:S) that were added to the atlas metadata
# sheet - the author backed this data out again on the next revision.
#
# ClearContents() drops the cached values, prunes now-unused shared
# strings, and (for un-styled cells) removes the <c> elements entirely
# while tightening each row's `spans` + the sheet `dimension`. The one
# cell in this block that carries a style (R43) is preserved as an empty,
# styled cell rather than being removed outright.
#
# The O:R block is cleared in three separate row-range passes (mirroring
# how the edit was actually made) rather than one single ClearContents
# across the whole O1:R67 span.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1:R32").ClearContents()
$ws.Range("O33:R48").ClearContents()
$ws.Range("O49:R67").ClearContents()
$ws.Range("S1:S67").ClearContents()

# Re-point the view at the region the author was working in when they saved.
$ws.Activate()
$ws.Range("E43").Select()
